$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 2539.3125
$ws.Cells.Item(17, 10).Value = 2539.3125
$ws.Cells.Item(17, 12).Value = 7617.9375
$ws.Cells.Item(17, 14).Value = -7953.9375
$ws.Cells.Item(40, 8).Value = 5470.7646
$ws.Cells.Item(40, 10).Value = 3000.7
$ws.Cells.Item(40, 12).Value = 3000.7
$ws.Cells.Item(40, 14).Value = -3350.7
$ws.Cells.Item(43, 8).Value = 2265.4375
$ws.Cells.Item(43, 9).Value = 1873.5
$ws.Cells.Item(43, 11).Value = 1873.5
$ws.Cells.Item(43, 13).Value = -1804.5
$ws.Cells.Item(112, 8).Value = 5341.4067
$ws.Cells.Item(112, 10).Value = 6000.18
$ws.Cells.Item(112, 12).Value = 18000.54
$ws.Cells.Item(112, 14).Value = -20216.54
$ws.Cells.Item(113, 8).Value = 4391.421
$ws.Cells.Item(113, 9).Value = 4132
$ws.Cells.Item(113, 11).Value = 4132
$ws.Cells.Item(113, 13).Value = -878
$ws.Cells.Item(116, 8).Value = 269553.25
$ws.Cells.Item(116, 9).Value = 84425
$ws.Cells.Item(116, 11).Value = 84425
$ws.Cells.Item(116, 13).Value = -80983
$ws.Cells.Item(132, 8).Value = 56123.953
$ws.Cells.Item(132, 9).Value = 60046.074
$ws.Cells.Item(132, 11).Value = 180138.222
$ws.Cells.Item(132, 13).Value = -177608.222
$ws.Cells.Item(133, 8).Value = 69824
$ws.Cells.Item(133, 10).Value = 69824
$ws.Cells.Item(133, 12).Value = 69824
$ws.Cells.Item(133, 14).Value = -79944
$ws.Cells.Item(134, 8).Value = 84154
$ws.Cells.Item(134, 10).Value = 84154
$ws.Cells.Item(134, 12).Value = 84154
$ws.Cells.Item(134, 14).Value = -94294
$ws.Cells.Item(139, 8).Value = 94997.5
$ws.Cells.Item(139, 10).Value = 94997.5
$ws.Cells.Item(139, 12).Value = 94997.5
$ws.Cells.Item(139, 14).Value = -105277.5
$ws.Cells.Item(140, 8).Value = 89825.28999999999
$ws.Cells.Item(140, 10).Value = 89825.28999999999
$ws.Cells.Item(140, 12).Value = 89825.28999999999
$ws.Cells.Item(140, 14).Value = -100185.29

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1348.14
$ws.Cells.Item(2, 9).Value = 1464.697
$ws.Cells.Item(2, 10).Value = 1121.8823
$ws.Cells.Item(2, 11).Value = 1464.697
$ws.Cells.Item(2, 12).Value = 1121.8823
$ws.Cells.Item(2, 13).Value = -1351.697
$ws.Cells.Item(2, 14).Value = -1347.8823
$ws.Cells.Item(32, 8).Value = 7464056.5
$ws.Cells.Item(32, 9).Value = 7937457.5
$ws.Cells.Item(32, 11).Value = 7937457.5
$ws.Cells.Item(32, 13).Value = -7937170.5
$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(41, 8).Value = 3061.75
$ws.Cells.Item(41, 9).Value = 0
$ws.Cells.Item(41, 11).Value = 0
$ws.Cells.Item(45, 8).Value = 3127.4375
$ws.Cells.Item(45, 9).Value = 3612.1667
$ws.Cells.Item(45, 10).Value = 1673.25
$ws.Cells.Item(45, 11).Value = 3612.1667
$ws.Cells.Item(45, 12).Value = 1673.25
$ws.Cells.Item(45, 13).Value = -3235.1667
$ws.Cells.Item(45, 14).Value = -2427.25
$ws.Cells.Item(74, 8).Value = 2121665.5
$ws.Cells.Item(74, 9).Value = 2502566.8
$ws.Cells.Item(74, 11).Value = 2502566.8
$ws.Cells.Item(74, 13).Value = -2501692.8
$ws.Cells.Item(77, 8).Value = 2121665.5
$ws.Cells.Item(77, 9).Value = 2502566.8
$ws.Cells.Item(77, 11).Value = 12512834
$ws.Cells.Item(77, 13).Value = -12508466
$ws.Cells.Item(116, 8).Value = 1348.14
$ws.Cells.Item(116, 9).Value = 1464.697
$ws.Cells.Item(116, 10).Value = 1121.8823
$ws.Cells.Item(116, 11).Value = 1464.697
$ws.Cells.Item(116, 12).Value = 1121.8823
$ws.Cells.Item(116, 13).Value = 829.3030000000001
$ws.Cells.Item(116, 14).Value = -5709.8823
$ws.Cells.Item(132, 8).Value = 865624.7
$ws.Cells.Item(132, 9).Value = 934478.75
$ws.Cells.Item(132, 11).Value = 2803436.25
$ws.Cells.Item(132, 13).Value = -2800906.25
$ws.Range("N40").ClearContents()
$ws.Range("M41").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1348.14
$ws.Cells.Item(3, 9).Value = 1464.697
$ws.Cells.Item(3, 10).Value = 1121.8823
$ws.Cells.Item(3, 11).Value = 1464.697
$ws.Cells.Item(3, 12).Value = 1121.8823
$ws.Cells.Item(3, 13).Value = -1350.697
$ws.Cells.Item(3, 14).Value = -1349.8823
$ws.Cells.Item(134, 8).Value = 995666.8
$ws.Cells.Item(134, 9).Value = 1833869.4
$ws.Cells.Item(134, 11).Value = 5501608.199999999
$ws.Cells.Item(134, 13).Value = -5499073.199999999
$ws.Cells.Item(140, 8).Value = 99097.5
$ws.Cells.Item(140, 10).Value = 99097.5
$ws.Cells.Item(140, 12).Value = 99097.5
$ws.Cells.Item(140, 14).Value = -109457.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 119951.77
$ws.Cells.Item(31, 9).Value = 191257.94
$ws.Cells.Item(31, 11).Value = 191257.94
$ws.Cells.Item(31, 13).Value = -190962.94
$ws.Cells.Item(34, 8).Value = 119951.77
$ws.Cells.Item(34, 9).Value = 191257.94
$ws.Cells.Item(34, 11).Value = 191257.94
$ws.Cells.Item(34, 13).Value = -191055.94
$ws.Cells.Item(58, 8).Value = 590533.0600000001
$ws.Cells.Item(58, 9).Value = 950533.9399999999
$ws.Cells.Item(58, 11).Value = 950533.9399999999
$ws.Cells.Item(58, 13).Value = -950330.9399999999
$ws.Cells.Item(132, 8).Value = 10433230
$ws.Cells.Item(132, 9).Value = 17933.137
$ws.Cells.Item(132, 10).Value = 125001500
$ws.Cells.Item(132, 11).Value = 53799.41099999999
$ws.Cells.Item(132, 12).Value = 375004500
$ws.Cells.Item(132, 13).Value = -51269.41099999999
$ws.Cells.Item(132, 14).Value = -375009560
$ws.Cells.Item(134, 8).Value = 13205.966
$ws.Cells.Item(134, 9).Value = 18589.578
$ws.Cells.Item(134, 10).Value = 2977.1
$ws.Cells.Item(134, 11).Value = 55768.734
$ws.Cells.Item(134, 12).Value = 8931.299999999999
$ws.Cells.Item(134, 13).Value = -53233.734
$ws.Cells.Item(134, 14).Value = -14001.3
$ws.Cells.Item(136, 8).Value = 590533.0600000001
$ws.Cells.Item(136, 9).Value = 950533.9399999999
$ws.Cells.Item(136, 11).Value = 2851601.82
$ws.Cells.Item(136, 13).Value = -2849051.82

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(39, 8).Value = 6333.3335
$ws.Cells.Item(39, 10).Value = 6333.3335
$ws.Cells.Item(39, 12).Value = 19000.0005
$ws.Cells.Item(39, 14).Value = -19588.0005
$ws.Cells.Item(55, 8).Value = 1887.7778
$ws.Cells.Item(55, 9).Value = 995
$ws.Cells.Item(55, 11).Value = 2985
$ws.Cells.Item(55, 13).Value = -2808

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 8100.25
$ws.Cells.Item(70, 9).Value = 9134
$ws.Cells.Item(70, 11).Value = 9134
$ws.Cells.Item(70, 13).Value = -8864
$ws.Cells.Item(73, 8).Value = 8100.25
$ws.Cells.Item(73, 9).Value = 9134
$ws.Cells.Item(73, 11).Value = 9134
$ws.Cells.Item(73, 13).Value = -8198
$ws.Cells.Item(102, 8).Value = 43481424
$ws.Cells.Item(102, 9).Value = 71430344
$ws.Cells.Item(102, 10).Value = 5332.5557
$ws.Cells.Item(102, 11).Value = 71430344
$ws.Cells.Item(102, 12).Value = 5332.5557
$ws.Cells.Item(102, 13).Value = -71428722
$ws.Cells.Item(102, 14).Value = -8576.555700000001
$ws.Cells.Item(132, 8).Value = 1206410.2
$ws.Cells.Item(132, 9).Value = 1339962
$ws.Cells.Item(132, 11).Value = 4019886
$ws.Cells.Item(132, 13).Value = -4017356

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3084
$ws.Cells.Item(7, 9).Value = 2900.8
$ws.Cells.Item(7, 11).Value = 2900.8
$ws.Cells.Item(7, 13).Value = -2788.8
$ws.Cells.Item(55, 8).Value = 893.1177
$ws.Cells.Item(55, 9).Value = 277.2857
$ws.Cells.Item(55, 10).Value = 1324.2
$ws.Cells.Item(55, 11).Value = 277.2857
$ws.Cells.Item(55, 12).Value = 1324.2
$ws.Cells.Item(55, 13).Value = -104.2857
$ws.Cells.Item(55, 14).Value = -1670.2
$ws.Cells.Item(61, 8).Value = 3479.52
$ws.Cells.Item(61, 9).Value = 1542.75
$ws.Cells.Item(61, 11).Value = 1542.75
$ws.Cells.Item(61, 13).Value = -1340.75
$ws.Cells.Item(68, 8).Value = 4262.5386
$ws.Cells.Item(68, 9).Value = 2233.3333
$ws.Cells.Item(68, 11).Value = 2233.3333
$ws.Cells.Item(68, 13).Value = -1484.3333
$ws.Cells.Item(71, 8).Value = 4262.5386
$ws.Cells.Item(71, 9).Value = 2233.3333
$ws.Cells.Item(71, 11).Value = 11166.6665
$ws.Cells.Item(71, 13).Value = -7422.666499999999
$ws.Cells.Item(113, 8).Value = 3479.52
$ws.Cells.Item(113, 9).Value = 1542.75
$ws.Cells.Item(113, 11).Value = 1542.75
$ws.Cells.Item(113, 13).Value = 627.25
$ws.Cells.Item(126, 8).Value = 3084
$ws.Cells.Item(126, 9).Value = 2900.8
$ws.Cells.Item(126, 11).Value = 8702.400000000001
$ws.Cells.Item(126, 13).Value = -6232.400000000001
$ws.Cells.Item(132, 8).Value = 1194919.1
$ws.Cells.Item(132, 9).Value = 1332390.5
$ws.Cells.Item(132, 10).Value = 3500
$ws.Cells.Item(132, 11).Value = 3997171.5
$ws.Cells.Item(132, 12).Value = 10500
$ws.Cells.Item(132, 13).Value = -3994641.5
$ws.Cells.Item(132, 14).Value = -15560
$ws.Cells.Item(136, 8).Value = 95567.42999999999
$ws.Cells.Item(136, 9).Value = 7428.4546
$ws.Cells.Item(136, 10).Value = 418743.66
$ws.Cells.Item(136, 11).Value = 22285.3638
$ws.Cells.Item(136, 12).Value = 1256230.98
$ws.Cells.Item(136, 13).Value = -19735.3638
$ws.Cells.Item(136, 14).Value = -1261330.98

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 6500.5
$ws.Cells.Item(96, 9).Value = 7351.5
$ws.Cells.Item(96, 11).Value = 7351.5
$ws.Cells.Item(96, 13).Value = -5978.5
$ws.Cells.Item(107, 8).Value = 2802.5625
$ws.Cells.Item(107, 9).Value = 1561.7084
$ws.Cells.Item(107, 11).Value = 4685.1252
$ws.Cells.Item(107, 13).Value = -2765.1252
$ws.Cells.Item(132, 8).Value = 16777684
$ws.Cells.Item(132, 9).Value = 50310560
$ws.Cells.Item(132, 10).Value = 11248.5
$ws.Cells.Item(132, 11).Value = 150931680
$ws.Cells.Item(132, 12).Value = 33745.5
$ws.Cells.Item(132, 13).Value = -150929150
$ws.Cells.Item(132, 14).Value = -38805.5
$ws.Cells.Item(136, 8).Value = 6236007
$ws.Cells.Item(136, 9).Value = 7451592
$ws.Cells.Item(136, 10).Value = 36525.2
$ws.Cells.Item(136, 11).Value = 22354776
$ws.Cells.Item(136, 12).Value = 109575.6
$ws.Cells.Item(136, 13).Value = -22352226
$ws.Cells.Item(136, 14).Value = -114675.6

Write-Host "Edit complete"